$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A9 holds a date-like string ("12/03/2025") that must stay plain text,
# matching the existing Date column cells (A2:A8) rather than being
# auto-converted into a date serial number by Excel's type inference.
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "12/03/2025"
$ws.Range("A9").ClearFormats()

$ws.Range("B9").Value = 14460.2
$ws.Range("C9").Value = 0.163896978737768
$ws.Range("D9").Value = 0.836103021262232
$ws.Range("E9").Value = -57.02
$ws.Range("F9").Value = -13.91
$ws.Range("G9").Value = -18376.54
$ws.Range("H9").Value = -60.32
$ws.Range("I9").Value = -432.87
$ws.Range("J9").Value = -15.44
